# "removed a weird annotation"
#
# Row 77 on the "Main script" sheet held a stray annotation row (the
# "Q5.0 / Nurses or Doctors" doctor-vs-nurse prompt, with its rich-text
# phoneme-tagged runs) that doesn't belong with the rest of the script.
# Delete that entire row; every row below it shifts up by one, and the
# sheet's filtered range / named ranges / selection need to follow.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main script")

# Remove the stray row entirely (shifts rows 78:157 up to 77:156).
$ws.Rows.Item(77).Delete()

# The filter/header named ranges covered through the old last row (151 /
# 157). After the deletion the data now ends one row earlier, so bring
# both defined names back in sync.
$wb.Names.Item("_xlnm._FilterDatabase").RefersTo = "='Main script'!`$E`$8:`$E`$150"
$wb.Names.Item("Z_E7F284B7_B5F7_4D04_B2BD_CA5521DD5FA3_.wvu.FilterData").RefersTo = "='Main script'!`$E`$1:`$E`$156"

# Re-apply the AutoFilter so its stored range matches the new extent too.
$ws.AutoFilterMode = $false
$ws.Range("E8:E150").AutoFilter()

# Reflect where the editor ended up after doing the deletion.
$ws.Range("A77").Select()
